# Auto-generated edit script
$d = $word.ActiveDocument

# --- Change 1: cached DATE field text in the document body ---
$d.Content.Find.Execute("December 22, 2021", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "January 10, 2022", 2) | Out-Null

# --- Change 2: footer "Vanderbilt Memory & Alzheimer's Center" line ---
# Drop the " · Vanderbilt Memory & Aging Project" tail and the
# trailing space after "ter".
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ftr.Range.Find.Execute("ter · Vanderbilt Memory & Aging Project", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ter", 2) | Out-Null

# --- Change 3: footer phone-number line ---
# Replace 615-347-6937 with 615-336-3388, split into its own run with an
# added complex-script size (szCs) of 9pt, matching the rest of the footer.
$find = $ftr.Range.Duplicate
$found = $find.Find.Execute("Nashville")
if ($found) {
    $insertPoint = $find.Duplicate
    $insertPoint.Collapse(1)
    $xmlFrag = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:pPr><w:tabs><w:tab w:val=`"center`" w:pos=`"4680`"/><w:tab w:val=`"right`" w:pos=`"9360`"/></w:tabs><w:spacing w:line=`"240`" w:lineRule=`"exact`"/><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:sz w:val=`"22`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:sz w:val=`"22`"/></w:rPr><w:t>1207 17</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:sz w:val=`"22`"/><w:vertAlign w:val=`"superscript`"/></w:rPr><w:t>th</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:sz w:val=`"22`"/></w:rPr><w:t xml:space=`"preserve`"> Avenue South, 2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:sz w:val=`"22`"/><w:vertAlign w:val=`"superscript`"/></w:rPr><w:t>nd</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:sz w:val=`"22`"/></w:rPr><w:t xml:space=`"preserve`"> floor, Suite 204 </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:b/><w:sz w:val=`"22`"/></w:rPr><w:t xml:space=`"preserve`">· </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:sz w:val=`"22`"/></w:rPr><w:t xml:space=`"preserve`">Nashville, TN 37212 </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:b/><w:sz w:val=`"22`"/></w:rPr><w:t>·</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:sz w:val=`"22`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Andalus`" w:eastAsia=`"FangSong`" w:hAnsi=`"Andalus`" w:cs=`"Andalus`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"18`"/></w:rPr><w:t>615-336-3388</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $insertPoint.InsertXML($xmlFrag)
}
